$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.44167
$ws.Range("I2").Value = 0.310972067878236
$ws.Range("J2").Value = 0.310972067878236
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 4.524359666666666
$ws.Range("N2").Value = 13.573079
$ws.Range("O2").Value = 0.9473211410445694
$ws.Range("P2").Value = 0.9473211410445695
$ws.Range("Q2").Value = 1.998273933976667
$ws.Range("R2").Value = 17.98446540579
$ws.Range("S2").Value = 0.2945904141753998
$ws.Range("T2").Value = 0.2945904141753999

# Row 3
$ws.Range("G3").Value = 0.44167
$ws.Range("I3").Value = 0.310972067878236
$ws.Range("J3").Value = 0.310972067878236
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.2515916666666667
$ws.Range("N3").Value = 0.754775
$ws.Range("O3").Value = 0.05267885895543045
$ws.Range("P3").Value = 0.05267885895543045
$ws.Range("Q3").Value = 0.1111204914166667
$ws.Range("R3").Value = 1.00008442275
$ws.Range("S3").Value = 0.01638165370283614
$ws.Range("T3").Value = 0.01638165370283614

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.8967930000000001
$ws.Range("H4").Value = 2.690379
$ws.Range("I4").Value = 0.6314161561091469
$ws.Range("J4").Value = 0.6314161561091469
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 4.524359666666666
$ws.Range("N4").Value = 13.573079
$ws.Range("O4").Value = 0.9473211410445694
$ws.Range("P4").Value = 0.9473211410445695
$ws.Range("Q4").Value = 4.057414078549
$ws.Range("R4").Value = 36.516726706941
$ws.Range("S4").Value = 0.5981538734792931
$ws.Range("T4").Value = 0.5981538734792932

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.8967930000000001
$ws.Range("H5").Value = 2.690379
$ws.Range("I5").Value = 0.6314161561091469
$ws.Range("J5").Value = 0.6314161561091469
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.2515916666666667
$ws.Range("N5").Value = 0.754775
$ws.Range("O5").Value = 0.05267885895543045
$ws.Range("P5").Value = 0.05267885895543045
$ws.Range("Q5").Value = 0.225625645525
$ws.Range("R5").Value = 2.030630809725
$ws.Range("S5").Value = 0.03326228262985381
$ws.Range("T5").Value = 0.03326228262985381

# Row 6
$ws.Range("G6").Value = 0.08182533333333333
$ws.Range("H6").Value = 0.245476
$ws.Range("I6").Value = 0.05761177601261715
$ws.Range("J6").Value = 0.05761177601261716
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 4.524359666666666
$ws.Range("N6").Value = 13.573079
$ws.Range("O6").Value = 0.9473211410445694
$ws.Range("P6").Value = 0.9473211410445695
$ws.Range("Q6").Value = 0.3702072378448888
$ws.Range("R6").Value = 3.331865140604
$ws.Range("S6").Value = 0.05457685338987664
$ws.Range("T6").Value = 0.05457685338987665

# Row 7
$ws.Range("G7").Value = 0.08182533333333333
$ws.Range("H7").Value = 0.245476
$ws.Range("I7").Value = 0.05761177601261715
$ws.Range("J7").Value = 0.05761177601261716
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.2515916666666667
$ws.Range("N7").Value = 0.754775
$ws.Range("O7").Value = 0.05267885895543045
$ws.Range("P7").Value = 0.05267885895543045
$ws.Range("Q7").Value = 0.02058657198888889
$ws.Range("R7").Value = 0.1852791479
$ws.Range("S7").Value = 0.00303492262274051
$ws.Range("T7").Value = 0.00303492262274051
